# Fix two typos in the "Suicides Among Veterans in the United States" deck.

$p = $ppt.ActivePresentation

# Slide 3 ("Questions and Motivation") - Content Placeholder 2, paragraph 3:
#   "...medical centers in the those high risk states?"
#                         -> "...medical centers in those high risk states?"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$run3 = $tr3.Paragraphs(3).Runs(1)
$run3.Text = "Which states have the most veterans" + [char]8217 + " death due to suicide and what is the total availability of Veteran Affairs (VA) medical centers in those high risk states?"

# Slide 22 ("Conclusion") - Content Placeholder 2, paragraph 5:
#   "Opening smaller non-profit clinics in rural parts..."
#                         -> "Opening smaller low cost/free clinics in rural parts..."
$s22 = $p.Slides.Item(22)
$tr22 = $s22.Shapes.Item(2).TextFrame.TextRange
$run22 = $tr22.Paragraphs(5).Runs(1)
$run22.Text = "Opening smaller low cost/free clinics in rural parts of those states that have high suicide rates in order to create accessibility for veterans to seek professional help."
